$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (want-to-go count) column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 2187
$wsExpo.Range("F3").Value = 904
$wsExpo.Range("F4").Value = 1607
$wsExpo.Range("F5").Value = 384

# Sheet "全部类型" (all types) - same underlying events, update column F accordingly
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 2187
$wsAll.Range("F5").Value = 904
$wsAll.Range("F6").Value = 1607
$wsAll.Range("F7").Value = 384
